$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("resp"), shifting resp/rt/word right
# into D/E/F, to make room for the new "name" column.
$ws.Columns("C:C").Insert()

# New header for the inserted column.
$ws.Range("C1").Value = "name"

# Row 2 is the instructions row: add its "name" and update resp/rt.
$ws.Range("C2").Value = "instr_resp"
$ws.Range("D2").Value = 13
$ws.Range("E2").Value = 1640.499999979511

# Rows 3-8 are stroop_response trials; fill "name" and updated resp/rt.
$ws.Range("C3").Value = "stroop_response"
$ws.Range("D3").Value = 119
$ws.Range("E3").Value = 168.7999999849126

$ws.Range("C4").Value = "stroop_response"
$ws.Range("D4").Value = 113
$ws.Range("E4").Value = 331.6999999806285

$ws.Range("C5").Value = "stroop_response"
$ws.Range("D5").Value = 119
$ws.Range("E5").Value = 249.8999999370426

$ws.Range("C6").Value = "stroop_response"
$ws.Range("D6").Value = 119
$ws.Range("E6").Value = 168.7999999849126

$ws.Range("C7").Value = "stroop_response"
$ws.Range("D7").Value = 113
$ws.Range("E7").Value = 331.6999999806285

$ws.Range("C8").Value = "stroop_response"
$ws.Range("D8").Value = 119
$ws.Range("E8").Value = 249.8999999370426
